$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.607064
$ws.Range("H2").Value = 4.821192
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.319612666666667
$ws.Range("N2").Value = 9.958838
$ws.Range("O2").Value = 0.1773188829476472
$ws.Range("P2").Value = 0.2112475282640173
$ws.Range("Q2").Value = 5.334830010544
$ws.Range("R2").Value = 48.013470094896
$ws.Range("S2").Value = 0.1773188829476472
$ws.Range("T2").Value = 0.2112475282640173
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.607064
$ws.Range("H3").Value = 4.821192
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.225632666666667
$ws.Range("N3").Value = 3.676898
$ws.Range("O3").Value = 0.0654678232613522
$ws.Range("P3").Value = 0.07799460280194422
$ws.Range("Q3").Value = 1.969670135824
$ws.Range("R3").Value = 17.727031222416
$ws.Range("S3").Value = 0.0654678232613522
$ws.Range("T3").Value = 0.07799460280194422
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.607064
$ws.Range("H4").Value = 4.821192
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.216913333333333
$ws.Range("N4").Value = 6.65074
$ws.Range("O4").Value = 0.1184176093210107
$ws.Range("P4").Value = 0.1410759353778654
$ws.Range("Q4").Value = 3.56272160912
$ws.Range("R4").Value = 32.06449448208
$ws.Range("S4").Value = 0.1184176093210107
$ws.Range("T4").Value = 0.1410759353778654
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.607064
$ws.Range("H5").Value = 4.821192
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.938529666666666
$ws.Range("N5").Value = 8.815588999999999
$ws.Range("O5").Value = 0.1569631310405457
$ws.Range("P5").Value = 0.1869968550991049
$ws.Range("Q5").Value = 4.722405240232
$ws.Range("R5").Value = 42.501647162088
$ws.Range("S5").Value = 0.1569631310405457
$ws.Range("T5").Value = 0.1869968550991049
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.607064
$ws.Range("H6").Value = 4.821192
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 9.020457499999999
$ws.Range("N6").Value = 18.040915
$ws.Range("O6").Value = 0.4818325534294442
$ws.Range("P6").Value = 0.3826850784570683
$ws.Range("Q6").Value = 14.49645251178
$ws.Range("R6").Value = 86.97871507068
$ws.Range("S6").Value = 0.4818325534294442
$ws.Range("T6").Value = 0.3826850784570683
